# Add a second column "Groepnaam (uniek)" with two values next to the
# existing Classificatiecode column, and move the active selection to A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Groepnaam (uniek)"
$ws.Range("B2").Value = "Inhang"
$ws.Range("B3").Value = "Kopschotten"

$ws.Range("A3").Select()
